$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header style (H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J for rows 2-51.
$data = @(
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(5, 5),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(5, 6),
    @(11, 11),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(8, 9),
    @(6, 6),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(5, 6),
    @(10, 11),
    @(8, 9),
    @(10, 10),
    @(7, 7),
    @(9, 9),
    @(6, 7),
    @(9, 10),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(4, 4),
    @(6, 6),
    @(7, 7)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
